# Add two new columns, I ("I0") and J ("IF"), to the sheet — mirroring the
# existing H ("IP") column: header in row 1 (same style as the other
# headers) and numeric values in rows 2-45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) -------------------------------------------------
# Copy the header formatting from H1 (bold, centered, bordered) onto the
# two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (2-45) --------------------------------------------------
$data = @{
    2  = @(1, 2)
    3  = @(1, 1)
    4  = @(1, 1)
    5  = @(1, 2)
    6  = @(3, 5)
    7  = @(1, 2)
    8  = @(5, 6)
    9  = @(2, 4)
    10 = @(3, 4)
    11 = @(9, 9)
    12 = @(6, 7)
    13 = @(8, 9)
    14 = @(9, 9)
    15 = @(1, 2)
    16 = @(7, 8)
    17 = @(6, 7)
    18 = @(6, 7)
    19 = @(7, 8)
    20 = @(7, 8)
    21 = @(7, 7)
    22 = @(7, 8)
    23 = @(7, 8)
    24 = @(6, 8)
    25 = @(9, 9)
    26 = @(6, 7)
    27 = @(6, 7)
    28 = @(6, 8)
    29 = @(6, 7)
    30 = @(6, 7)
    31 = @(9, 9)
    32 = @(3, 5)
    33 = @(8, 9)
    34 = @(6, 7)
    35 = @(7, 8)
    36 = @(8, 9)
    37 = @(7, 8)
    38 = @(4, 5)
    39 = @(4, 6)
    40 = @(5, 6)
    41 = @(5, 6)
    42 = @(5, 6)
    43 = @(8, 9)
    44 = @(3, 4)
    45 = @(1, 2)
}

foreach ($row in 2..45) {
    $values = $data[$row]
    $ws.Range("I$row").Value = $values[0]
    $ws.Range("J$row").Value = $values[1]
}

Write-Output "I0/IF columns added"
